$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "total" header label for the procura/não procuraram columns was
# mislabeled as the stray pandas artifact "unnamed: 1_level_1"; fix it to
# match the "total" label used in B1.
$ws.Range("B2").Value = "total"

# Two sub-header rows ("situação do domicílio" and "grandes regiões e
# unidades da federação") were left as empty placeholder rows with no data
# underneath them. Remove them so the data rows that follow each one shift
# up directly beneath the correct top-level row ("brasil" / "rural").
# Delete the lower row first so the upper row's index ("5") stays valid.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()
